$wb = $excel.ActiveWorkbook

# --- ShipDB ---------------------------------------------------------------
$wsShip = $wb.Worksheets.Item("ShipDB")

# S0011 / A-Wing qty 1 -> 2
$wsShip.Range("C12").Value = 2

$wsShip.Range("A21").Value = 'S0020'
$wsShip.Range("B21").Value = 'K-Wing'
$wsShip.Range("C21").Value = 1

$wsShip.Range("A22").Value = 'S0021'
$wsShip.Range("B22").Value = 'Scurrg H-6 Bomber'
$wsShip.Range("C22").Value = 1

# --- PilotDB ----------------------------------------------------------------
$wsPilot = $wb.Worksheets.Item("PilotDB")

# Green Squadron Pilot / Prototype Pilot (A-Wing) qty 1 -> 2
$wsPilot.Range("M51").Value = 2
$wsPilot.Range("M52").Value = 2

# Row 97: Tycho Celchu
$wsPilot.Range("A97").Value = 'P0096'
$wsPilot.Range("B97").Value = 'Tycho Celchu'
$wsPilot.Range("C97").Value = 'You may perform actions even while you have stress tokens.'
$wsPilot.Range("D97").Value = 26
$wsPilot.Range("E97").Value = 'Y'
$wsPilot.Range("F97").Value = 'A-Wing'
$wsPilot.Range("G97").Value = 'Rebel'
$wsPilot.Range("H97").Value = 8
$wsPilot.Range("I97").Value = 2
$wsPilot.Range("J97").Value = 3
$wsPilot.Range("K97").Value = 2
$wsPilot.Range("L97").Value = 2
$wsPilot.Range("M97").Value = 1

# Row 98: Arvel Crynyd
$wsPilot.Range("A98").Value = 'P0097'
$wsPilot.Range("B98").Value = 'Arvel Crynyd'
$wsPilot.Range("C98").Value = 'You may declare an enemy ship inside your firing arc that you are touching as the target of your attack.'
$wsPilot.Range("D98").Value = 23
$wsPilot.Range("E98").Value = 'Y'
$wsPilot.Range("F98").Value = 'A-Wing'
$wsPilot.Range("G98").Value = 'Rebel'
$wsPilot.Range("H98").Value = 6
$wsPilot.Range("I98").Value = 2
$wsPilot.Range("J98").Value = 3
$wsPilot.Range("K98").Value = 2
$wsPilot.Range("L98").Value = 2
$wsPilot.Range("M98").Value = 1

# Row 99: Miranda Doni
$wsPilot.Range("A99").Value = 'P0098'
$wsPilot.Range("B99").Value = 'Miranda Doni'
$wsPilot.Range("C99").Value = 'Once per round when attacking, you may either spend 1 shield to roll 1 additional attack die OR roll 1 fewer attack die to recover 1 shield.'
$wsPilot.Range("D99").Value = 29
$wsPilot.Range("E99").Value = 'Y'
$wsPilot.Range("F99").Value = 'K-Wing'
$wsPilot.Range("G99").Value = 'Rebel'
$wsPilot.Range("H99").Value = 8
$wsPilot.Range("I99").Value = 2
$wsPilot.Range("J99").Value = 1
$wsPilot.Range("K99").Value = 5
$wsPilot.Range("L99").Value = 4
$wsPilot.Range("M99").Value = 1

# Row 100: Esege Tuketu
$wsPilot.Range("A100").Value = 'P0099'
$wsPilot.Range("B100").Value = 'Esege Tuketu'
$wsPilot.Range("C100").Value = 'When another friendly ship at Range 1-2 is attacking, it may treat your focus tokens as its own.'
$wsPilot.Range("D100").Value = 28
$wsPilot.Range("E100").Value = 'Y'
$wsPilot.Range("F100").Value = 'K-Wing'
$wsPilot.Range("G100").Value = 'Rebel'
$wsPilot.Range("H100").Value = 6
$wsPilot.Range("I100").Value = 2
$wsPilot.Range("J100").Value = 1
$wsPilot.Range("K100").Value = 5
$wsPilot.Range("L100").Value = 4
$wsPilot.Range("M100").Value = 1

# Row 101: Guardian Squadron Pilot
$wsPilot.Range("A101").Value = 'P0100'
$wsPilot.Range("B101").Value = 'Guardian Squadron Pilot'
$wsPilot.Range("C101").Value = '-'
$wsPilot.Range("D101").Value = 25
$wsPilot.Range("E101").Value = 'N'
$wsPilot.Range("F101").Value = 'K-Wing'
$wsPilot.Range("G101").Value = 'Rebel'
$wsPilot.Range("H101").Value = 4
$wsPilot.Range("I101").Value = 2
$wsPilot.Range("J101").Value = 1
$wsPilot.Range("K101").Value = 5
$wsPilot.Range("L101").Value = 4
$wsPilot.Range("M101").Value = 1

# Row 102: Warden Squadron Pilot
$wsPilot.Range("A102").Value = 'P0101'
$wsPilot.Range("B102").Value = 'Warden Squadron Pilot'
$wsPilot.Range("C102").Value = '-'
$wsPilot.Range("D102").Value = 23
$wsPilot.Range("E102").Value = 'N'
$wsPilot.Range("F102").Value = 'K-Wing'
$wsPilot.Range("G102").Value = 'Rebel'
$wsPilot.Range("H102").Value = 2
$wsPilot.Range("I102").Value = 2
$wsPilot.Range("J102").Value = 1
$wsPilot.Range("K102").Value = 5
$wsPilot.Range("L102").Value = 4
$wsPilot.Range("M102").Value = 1

# Row 103: Captain Nym
$wsPilot.Range("A103").Value = 'P0102'
$wsPilot.Range("B103").Value = 'Captain Nym'
$wsPilot.Range("C103").Value = 'Once per round, you may prevent a friendly bomb from detonating.'
$wsPilot.Range("D103").Value = 30
$wsPilot.Range("E103").Value = 'Y'
$wsPilot.Range("F103").Value = 'Scurrg H-6 Bomber'
$wsPilot.Range("G103").Value = 'Rebel'
$wsPilot.Range("H103").Value = 8
$wsPilot.Range("I103").Value = 3
$wsPilot.Range("J103").Value = 1
$wsPilot.Range("K103").Value = 5
$wsPilot.Range("L103").Value = 5
$wsPilot.Range("M103").Value = 1

# Row 104: Captain Nym
$wsPilot.Range("A104").Value = 'P0103'
$wsPilot.Range("B104").Value = 'Captain Nym'
$wsPilot.Range("C104").Value = 'You may ignore friendly bombs. When a friendly ship is defending, if the attacker measures range through a friendly bomb token, the defender may add 1 EVADE result.'
$wsPilot.Range("D104").Value = 30
$wsPilot.Range("E104").Value = 'Y'
$wsPilot.Range("F104").Value = 'Scurrg H-6 Bomber'
$wsPilot.Range("G104").Value = 'Scum'
$wsPilot.Range("H104").Value = 8
$wsPilot.Range("I104").Value = 3
$wsPilot.Range("J104").Value = 1
$wsPilot.Range("K104").Value = 5
$wsPilot.Range("L104").Value = 5
$wsPilot.Range("M104").Value = 1

# Row 105: Sol Sixxa
$wsPilot.Range("A105").Value = 'P0104'
$wsPilot.Range("B105").Value = 'Sol Sixxa'
$wsPilot.Range("C105").Value = 'When dropping a bomb, you may use the [TURN 1] template instead of the [FORWARD 1] template.'
$wsPilot.Range("D105").Value = 28
$wsPilot.Range("E105").Value = 'Y'
$wsPilot.Range("F105").Value = 'Scurrg H-6 Bomber'
$wsPilot.Range("G105").Value = 'Scum'
$wsPilot.Range("H105").Value = 6
$wsPilot.Range("I105").Value = 3
$wsPilot.Range("J105").Value = 1
$wsPilot.Range("K105").Value = 5
$wsPilot.Range("L105").Value = 5
$wsPilot.Range("M105").Value = 1

# Row 106: Lok Revenant
$wsPilot.Range("A106").Value = 'P0105'
$wsPilot.Range("B106").Value = 'Lok Revenant'
$wsPilot.Range("C106").Value = '-'
$wsPilot.Range("D106").Value = 26
$wsPilot.Range("E106").Value = 'N'
$wsPilot.Range("F106").Value = 'Scurrg H-6 Bomber'
$wsPilot.Range("G106").Value = 'Scum'
$wsPilot.Range("H106").Value = 3
$wsPilot.Range("I106").Value = 3
$wsPilot.Range("J106").Value = 1
$wsPilot.Range("K106").Value = 5
$wsPilot.Range("L106").Value = 5
$wsPilot.Range("M106").Value = 1

# Row 107: Karthakk Pirate
$wsPilot.Range("A107").Value = 'P0106'
$wsPilot.Range("B107").Value = 'Karthakk Pirate'
$wsPilot.Range("C107").Value = '-'
$wsPilot.Range("D107").Value = 24
$wsPilot.Range("E107").Value = 'N'
$wsPilot.Range("F107").Value = 'Scurrg H-6 Bomber'
$wsPilot.Range("G107").Value = 'Scum'
$wsPilot.Range("H107").Value = 1
$wsPilot.Range("I107").Value = 3
$wsPilot.Range("J107").Value = 1
$wsPilot.Range("K107").Value = 5
$wsPilot.Range("L107").Value = 5
$wsPilot.Range("M107").Value = 1

# --- CardDB -----------------------------------------------------------------
$wsCard = $wb.Worksheets.Item("CardDB")

# Seismic Torpedo / Assault Missiles qty 1 -> 2
$wsCard.Range("H68").Value = 2
$wsCard.Range("H83").Value = 2

# Row 113: Concussion Missiles
$wsCard.Range("A113").Value = 'U0112'
$wsCard.Range("B113").Value = 'Concussion Missiles'
$wsCard.Range("C113").Value = 'FP: 4, RNG: 2-3 ATTACK (TARGET LOCK): Spend the target lock and discard this card to perfrorm this attack. You may change one of your blank results for a HIT result.'
$wsCard.Range("D113").Value = 4
$wsCard.Range("E113").Value = 'N'
$wsCard.Range("F113").Value = 'N'
$wsCard.Range("G113").Value = 'Missile'
$wsCard.Range("H113").Value = 1

# Row 114: Cluster Missiles
$wsCard.Range("A114").Value = 'U0113'
$wsCard.Range("B114").Value = 'Cluster Missiles'
$wsCard.Range("C114").Value = 'FP: 3, RNG: 1-2 ATTACK (TARGET LOCK): Spend your target lock and discard this card to perform this attack TWICE.'
$wsCard.Range("D114").Value = 4
$wsCard.Range("E114").Value = 'N'
$wsCard.Range("F114").Value = 'N'
$wsCard.Range("G114").Value = 'Missile'
$wsCard.Range("H114").Value = 1

# Row 115: Push The Limit
$wsCard.Range("A115").Value = 'U0114'
$wsCard.Range("B115").Value = 'Push The Limit'
$wsCard.Range("C115").Value = 'Once per round, after you perform an action, you may perform 1 free action shown in your action bar. Then receive 1 stress token.'
$wsCard.Range("D115").Value = 3
$wsCard.Range("E115").Value = 'N'
$wsCard.Range("F115").Value = 'N'
$wsCard.Range("G115").Value = 'Elite Pilot Skill'
$wsCard.Range("H115").Value = 1

# Row 116: Deadeye
$wsCard.Range("A116").Value = 'U0115'
$wsCard.Range("B116").Value = 'Deadeye'
$wsCard.Range("C116").Value = 'You may treat the ''ATTACK (TARGET LOCK)'' header as ''ATTACK (FOCUS)''. When an attack instructs you to spend a target lock, you may spend a focus token instead.'
$wsCard.Range("D116").Value = 1
$wsCard.Range("E116").Value = 'N'
$wsCard.Range("F116").Value = 'N'
$wsCard.Range("G116").Value = 'Elite Pilot Skill'
$wsCard.Range("H116").Value = 1

# Row 117: Plasma Torpedoes
$wsCard.Range("A117").Value = 'U0116'
$wsCard.Range("B117").Value = 'Plasma Torpedoes'
$wsCard.Range("C117").Value = 'FP: 4, RNG: 2-3 ATTACK (TARGET LOCK): Spend your target lock and discard this card ot perform this attack. If this attack hits, After dealing damage, remove 1 shield token from the defender.'
$wsCard.Range("D117").Value = 3
$wsCard.Range("E117").Value = 'N'
$wsCard.Range("F117").Value = 'N'
$wsCard.Range("G117").Value = 'Torpedo'
$wsCard.Range("H117").Value = 1

# Row 118: Twin Laser Turret
$wsCard.Range("A118").Value = 'U0117'
$wsCard.Range("B118").Value = 'Twin Laser Turret'
$wsCard.Range("C118").Value = 'FP: 3, RNG: 2-3 ATTACK: Perform this attack TWICE (even against a ship outside your firing arc). Each time this attack hits, the defender suffers 1 damage. Then cancel ALL dice results.'
$wsCard.Range("D118").Value = 6
$wsCard.Range("E118").Value = 'N'
$wsCard.Range("F118").Value = 'N'
$wsCard.Range("G118").Value = 'Turret'
$wsCard.Range("H118").Value = 2

# Row 119: Adv. Homing Missiles
$wsCard.Range("A119").Value = 'U0118'
$wsCard.Range("B119").Value = 'Adv. Homing Missiles'
$wsCard.Range("C119").Value = 'FP: 3, RNG 2 ATTACK (TARGET LOCK): Discard this card to perfrom this attack. If this attack hits, deal 1 face-up damage card to the defender. Then cancel ALL dice results.'
$wsCard.Range("D119").Value = 3
$wsCard.Range("E119").Value = 'N'
$wsCard.Range("F119").Value = 'N'
$wsCard.Range("G119").Value = 'Missile'
$wsCard.Range("H119").Value = 1

# Row 120: Bombardier
$wsCard.Range("A120").Value = 'U0119'
$wsCard.Range("B120").Value = 'Bombardier'
$wsCard.Range("C120").Value = 'When dropping a bomb, you may use the [FORWARD 2] template instead of the [FORWARD 1] template.'
$wsCard.Range("D120").Value = 1
$wsCard.Range("E120").Value = 'N'
$wsCard.Range("F120").Value = 'N'
$wsCard.Range("G120").Value = 'Crew'
$wsCard.Range("H120").Value = 1

# Row 121: Conner Net
$wsCard.Range("A121").Value = 'U0120'
$wsCard.Range("B121").Value = 'Conner Net'
$wsCard.Range("C121").Value = 'ACTION: Discard this card to DROP 1 Conner net token. When a ship''s base or maneuver template overlaps this token, this token DETONATES.'
$wsCard.Range("D121").Value = 4
$wsCard.Range("E121").Value = 'N'
$wsCard.Range("F121").Value = 'N'
$wsCard.Range("G121").Value = 'Bomb'
$wsCard.Range("H121").Value = 1

# Row 122: Extra Munitions
$wsCard.Range("A122").Value = 'U0121'
$wsCard.Range("B122").Value = 'Extra Munitions'
$wsCard.Range("C122").Value = 'When you equip this card, place 1 ordnance token on each equipped TORPEDO, MISSILE and BOMB Upgrade card. When you are instructed to discard an upgrade card, you may discard 1 ordnance token on that card instead.'
$wsCard.Range("D122").Value = 2
$wsCard.Range("E122").Value = 'N'
$wsCard.Range("F122").Value = 'Y'
$wsCard.Range("G122").Value = 'Torpedo'
$wsCard.Range("H122").Value = 1

# Row 123: Ion Bombs
$wsCard.Range("A123").Value = 'U0122'
$wsCard.Range("B123").Value = 'Ion Bombs'
$wsCard.Range("C123").Value = 'When you reveal your maneuver dial, you may discard this card to DROP 1 ion bomb token. This token detonates at the end of the activation phase.'
$wsCard.Range("D123").Value = 2
$wsCard.Range("E123").Value = 'N'
$wsCard.Range("F123").Value = 'N'
$wsCard.Range("G123").Value = 'Bomb'
$wsCard.Range("H123").Value = 1

# Row 124: Advanced SLAM
$wsCard.Range("A124").Value = 'U0123'
$wsCard.Range("B124").Value = 'Advanced SLAM'
$wsCard.Range("C124").Value = 'After performing a SLAM action, if you did not overlap an obstacle or another ship, you may perform a free action'
$wsCard.Range("D124").Value = 2
$wsCard.Range("E124").Value = 'N'
$wsCard.Range("F124").Value = 'N'
$wsCard.Range("G124").Value = 'Modification'
$wsCard.Range("H124").Value = 1

# Row 125: Lightning Reflexes
$wsCard.Range("A125").Value = 'U0124'
$wsCard.Range("B125").Value = 'Lightning Reflexes'
$wsCard.Range("C125").Value = 'SMALL SHIP ONLY. After you execute a white or green maneuveron your dial, you may discard this card to rotate your ship 180 degrees. Then receive 1 stress token AFTER the ''check pilot stress'' step.'
$wsCard.Range("D125").Value = 1
$wsCard.Range("E125").Value = 'N'
$wsCard.Range("F125").Value = 'N'
$wsCard.Range("G125").Value = 'Elite Pilot Skill'
$wsCard.Range("H125").Value = 1

# Row 126: Minefield Mapper
$wsCard.Range("A126").Value = 'U0125'
$wsCard.Range("B126").Value = 'Minefield Mapper'
$wsCard.Range("C126").Value = 'During setup, after the ''place forces'' step, you may discard any number of equipped BOMB upgrade cards. Place all corresponding bomb tokens in the play area beyond Range 3 of enemy ships.'
$wsCard.Range("D126").Value = 0
$wsCard.Range("E126").Value = 'N'
$wsCard.Range("F126").Value = 'N'
$wsCard.Range("G126").Value = 'System'
$wsCard.Range("H126").Value = 1

# Row 127: Synced Turret
$wsCard.Range("A127").Value = 'U0126'
$wsCard.Range("B127").Value = 'Synced Turret'
$wsCard.Range("C127").Value = 'FP: 3, RNG 1-2 ATTACK (TARGET LOCK): Attack 1 ship (even a ship outside yoru firing arc). If the defender is inside your primary firing arc, you may reroll a number of attack dice up to your primary weapoon value.'
$wsCard.Range("D127").Value = 4
$wsCard.Range("E127").Value = 'N'
$wsCard.Range("F127").Value = 'N'
$wsCard.Range("G127").Value = 'Turret'
$wsCard.Range("H127").Value = 1

# Row 128: Cruise Missiles
$wsCard.Range("A128").Value = 'U0127'
$wsCard.Range("B128").Value = 'Cruise Missiles'
$wsCard.Range("C128").Value = 'FP: 1, RNG: 2-3 ATTACK (TARGET LOCK): Discard this card to perform this attack. You may roll additional attack dice equal to the speed of the maneuver you executed this round, to a maximum of 4 additional dice.'
$wsCard.Range("D128").Value = 3
$wsCard.Range("E128").Value = 'N'
$wsCard.Range("F128").Value = 'N'
$wsCard.Range("G128").Value = 'Missile'
$wsCard.Range("H128").Value = 2

# Row 129: Cad Bane
$wsCard.Range("A129").Value = 'U0128'
$wsCard.Range("B129").Value = 'Cad Bane'
$wsCard.Range("C129").Value = 'SCUM ONLY. Your upgrade bar gains the BOMB icon. Once per round, when an enemy ship rolls attack dice due to a friendly bomb detonating, you may choose any number of FOCUS and BLANK results. It must reroll those results.'
$wsCard.Range("D129").Value = 2
$wsCard.Range("E129").Value = 'Y'
$wsCard.Range("F129").Value = 'N'
$wsCard.Range("G129").Value = 'Crew'
$wsCard.Range("H129").Value = 1

# Row 130: Bomblet Generator
$wsCard.Range("A130").Value = 'U0129'
$wsCard.Range("B130").Value = 'Bomblet Generator'
$wsCard.Range("C130").Value = 'When you reveal your maneuver, you may DROP 1 bomblet token. This token detonates at the end of the activation phase. (REQUIRES 2xBOMB SLOTS)'
$wsCard.Range("D130").Value = 3
$wsCard.Range("E130").Value = 'Y'
$wsCard.Range("F130").Value = 'N'
$wsCard.Range("G130").Value = 'Bomb'
$wsCard.Range("H130").Value = 1

# Row 131: R4-E1
$wsCard.Range("A131").Value = 'U0130'
$wsCard.Range("B131").Value = 'R4-E1'
$wsCard.Range("C131").Value = 'You can perform actions on you TORPEDO and BOMB Upgrade cards even if you are stressed. After you perform an action in this way, you may discard this card to remove 1 stress token from your ship.'
$wsCard.Range("D131").Value = 1
$wsCard.Range("E131").Value = 'Y'
$wsCard.Range("F131").Value = 'N'
$wsCard.Range("G131").Value = 'Salvaged Astromech'
$wsCard.Range("H131").Value = 1

# Row 132: Havoc
$wsCard.Range("A132").Value = 'U0131'
$wsCard.Range("B132").Value = 'Havoc'
$wsCard.Range("C132").Value = 'SCURRG H-6 BOMBER ONLY. Your upgrade bar gains the SYSTEM and SALVAGED ASTROMECH upgrade icons, and loses the CREW  upgrade icon. You cannot equip non-unique SALVAGED ASTROMECH Upgrade cards.'
$wsCard.Range("D132").Value = 0
$wsCard.Range("E132").Value = 'Y'
$wsCard.Range("F132").Value = 'N'
$wsCard.Range("G132").Value = 'Title'
$wsCard.Range("H132").Value = 1

# --- selection bookkeeping (matches sheetView changes in the diff) ---------
$wsShip.Range("B25").Select()
$wsPilot.Range("B109").Select()
$wsCard.Range("A110:A132").Select()
$wsCard.Select()

